$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 44383.5343287037

# Row 2
$ws.Cells.Item(2, 1).Value = $newDate

# Row 3
$ws.Cells.Item(3, 1).Value = $newDate
$ws.Cells.Item(3, 4).Value = 896.05924946
$ws.Cells.Item(3, 5).Value = 123.63124776
$ws.Cells.Item(3, 6).Value = 772.4280017
$ws.Cells.Item(3, 8).Value = 896.0599999999999

# Row 4
$ws.Cells.Item(4, 1).Value = $newDate
$ws.Cells.Item(4, 7).Value = 34312.7
$ws.Cells.Item(4, 8).Value = 528.47

# Row 5
$ws.Cells.Item(5, 1).Value = $newDate
$ws.Cells.Item(5, 4).Value = 22.0802
$ws.Cells.Item(5, 6).Value = 21.843
$ws.Cells.Item(5, 7).Value = 13.3314
$ws.Cells.Item(5, 8).Value = 294.36

# Row 6
$ws.Cells.Item(6, 1).Value = $newDate
$ws.Cells.Item(6, 7).Value = 0.9287

# Row 7
$ws.Cells.Item(7, 1).Value = $newDate
$ws.Cells.Item(7, 7).Value = 2325.29

# Row 8
$ws.Cells.Item(8, 1).Value = $newDate
